$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new weekly records ("Primera"/"Segunda" for Repollo,
# Crespo record @ Femacal de La Calera, dated 45041) that need to be inserted
# right before the existing row 905, pushing rows 905-952 down to 907-954.
$ws.Rows.Item(905).Resize(2).Insert()

# Populate the two newly-inserted rows (905 = "Primera", 906 = "Segunda")
# with the same constant/categorical values used throughout this sheet.
$ws.Range("A905:A906").Value = 3
$ws.Range("B905:B906").Value = "Femacal de La Calera"
$ws.Range("C905:C906").Value = "Coquimbo"
$ws.Range("E905:E906").Value = 5
$ws.Range("F905:F906").Value = 100112006
$ws.Range("G905:G906").Value = "Repollo"
$ws.Range("H905:H906").Value = "Crespo record"
$ws.Range("N905:N906").Value = "`$/unidad"
$ws.Range("O905:O906").Value = "Provincia de Quillota"
$ws.Range("Q905:Q906").Value = 1
$ws.Range("R905:R906").Value = "Hortaliza"

$ws.Range("I905").Value = "Primera"
$ws.Range("I906").Value = "Segunda"

# Row 905 ("Primera") data
$ws.Range("D905").Value = 45041
$ws.Range("J905").Value = 3200
$ws.Range("K905").Value = 1200
$ws.Range("L905").Value = 1300
$ws.Range("M905").Value = 1250
$ws.Range("P905").Value = 1250

# Row 906 ("Segunda") data
$ws.Range("D906").Value = 45041
$ws.Range("J906").Value = 1500
$ws.Range("K906").Value = 1000
$ws.Range("L906").Value = 1000
$ws.Range("M906").Value = 1000
$ws.Range("P906").Value = 1000

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D905:D906").NumberFormat = $ws.Range("D907").NumberFormat
